$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.542.33"
$ws.Range("E2").Value = "  -0.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.754.91"
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.41"
$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4569"
$ws.Range("E7").Value = "  +1.86%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3550"
$ws.Range("E8").Value = "  -1.76%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07472"
$ws.Range("E9").Value = "  -0.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.48"
$ws.Range("E10").Value = "  -1.55%  "

$ws.Range("E11").Value = "  -1.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.81"
$ws.Range("E13").Value = "  +0.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.009"
$ws.Range("E14").Value = "  -0.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.159"
$ws.Range("E15").Value = "  -0.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.754.46"
$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("E18").Value = "  -1.03%  "

$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9996"
$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.09"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.738"
$ws.Range("E22").Value = "  -2.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.592.37"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.22"
$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("E25").Value = "  -1.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.40"
$ws.Range("E26").Value = "  +2.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.16"
$ws.Range("E27").Value = "  -1.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.956.92"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.134"
$ws.Range("E29").Value = "  +0.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.58"
$ws.Range("E30").Value = "  +0.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.083"
$ws.Range("E31").Value = "  -0.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09237"
$ws.Range("E32").Value = "  +2.25%  "

$ws.Range("E33").Value = "  +0.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.533"
$ws.Range("E34").Value = "  -0.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.74"
$ws.Range("E35").Value = "  -2.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02285"
$ws.Range("E36").Value = "  -1.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2092"
$ws.Range("E37").Value = "  +0.21%  "

$ws.Range("E38").Value = "  +0.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6290"
$ws.Range("E39").Value = "  -1.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.924"
$ws.Range("E40").Value = "  -0.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.183"
$ws.Range("E41").Value = "  -1.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.388"
$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.806"
$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.16"
$ws.Range("E44").Value = "  -1.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.717"
$ws.Range("E45").Value = "  +0.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5866"
$ws.Range("E46").Value = "  -0.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.28"
$ws.Range("E47").Value = "  +0.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.939"
$ws.Range("E48").Value = "  -1.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06895"
$ws.Range("E49").Value = "  +0.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.130"
$ws.Range("E50").Value = "  -2.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.30"
$ws.Range("E51").Value = "  -0.20%  "
